# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" message text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.0 = 52926.7 pesos`n✅ 52926.7 pesos = 12.94 = 970.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update N10, O10, N12, O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 76.899
$ws2.Range("O10").Value = 4070.01
$ws2.Range("N12").Value = 4089.99
$ws2.Range("O12").Value = 75
